$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 184.0626906666667
$ws.Cells.Item(2, 8).Value = 552.188072
$ws.Cells.Item(2, 9).Value = 0.6510505751503485
$ws.Cells.Item(2, 10).Value = 0.6510505751503486
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 2.614605666666666
$ws.Cells.Item(2, 14).Value = 7.843817
$ws.Cells.Item(2, 15).Value = 0.07238665237615237
$ws.Cells.Item(2, 16).Value = 0.07238665237615237
$ws.Cells.Item(2, 17).Value = 481.2513540389804
$ws.Cells.Item(2, 18).Value = 4331.262186350824
$ws.Cells.Item(2, 19).Value = 0.04712737166270234
$ws.Cells.Item(2, 20).Value = 0.04712737166270235

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 184.0626906666667
$ws.Cells.Item(3, 8).Value = 552.188072
$ws.Cells.Item(3, 9).Value = 0.6510505751503485
$ws.Cells.Item(3, 10).Value = 0.6510505751503486
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 14.70158366666666
$ws.Cells.Item(3, 14).Value = 44.10475099999999
$ws.Cells.Item(3, 15).Value = 0.4070206225838464
$ws.Cells.Item(3, 16).Value = 0.4070206225838464
$ws.Cells.Item(3, 17).Value = 2706.013046747785
$ws.Cells.Item(3, 18).Value = 24354.11742073007
$ws.Cells.Item(3, 19).Value = 0.2649910104312661
$ws.Cells.Item(3, 20).Value = 0.2649910104312662

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 184.0626906666667
$ws.Cells.Item(4, 8).Value = 552.188072
$ws.Cells.Item(4, 9).Value = 0.6510505751503485
$ws.Cells.Item(4, 10).Value = 0.6510505751503486
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 18.80380766666667
$ws.Cells.Item(4, 14).Value = 56.411423
$ws.Cells.Item(4, 15).Value = 0.5205927250400012
$ws.Cells.Item(4, 16).Value = 0.5205927250400012
$ws.Cells.Item(4, 17).Value = 3461.079433905162
$ws.Cells.Item(4, 18).Value = 31149.71490514646
$ws.Cells.Item(4, 19).Value = 0.33893219305638
$ws.Cells.Item(4, 20).Value = 0.3389321930563801

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 57.4434
$ws.Cells.Item(5, 8).Value = 172.3302
$ws.Cells.Item(5, 9).Value = 0.2031838091312023
$ws.Cells.Item(5, 10).Value = 0.2031838091312023
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 2.614605666666666
$ws.Cells.Item(5, 14).Value = 7.843817
$ws.Cells.Item(5, 15).Value = 0.07238665237615237
$ws.Cells.Item(5, 16).Value = 0.07238665237615237
$ws.Cells.Item(5, 17).Value = 150.1918391526
$ws.Cells.Item(5, 18).Value = 1351.7265523734
$ws.Cells.Item(5, 19).Value = 0.01470779576004283
$ws.Cells.Item(5, 20).Value = 0.01470779576004283

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 57.4434
$ws.Cells.Item(6, 8).Value = 172.3302
$ws.Cells.Item(6, 9).Value = 0.2031838091312023
$ws.Cells.Item(6, 10).Value = 0.2031838091312023
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 14.70158366666666
$ws.Cells.Item(6, 14).Value = 44.10475099999999
$ws.Cells.Item(6, 15).Value = 0.4070206225838464
$ws.Cells.Item(6, 16).Value = 0.4070206225838464
$ws.Cells.Item(6, 17).Value = 844.5089511977998
$ws.Cells.Item(6, 18).Value = 7600.580560780199
$ws.Cells.Item(6, 19).Value = 0.08270000049153937
$ws.Cells.Item(6, 20).Value = 0.08270000049153937

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 57.4434
$ws.Cells.Item(7, 8).Value = 172.3302
$ws.Cells.Item(7, 9).Value = 0.2031838091312023
$ws.Cells.Item(7, 10).Value = 0.2031838091312023
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 18.80380766666667
$ws.Cells.Item(7, 14).Value = 56.411423
$ws.Cells.Item(7, 15).Value = 0.5205927250400012
$ws.Cells.Item(7, 16).Value = 0.5205927250400012
$ws.Cells.Item(7, 17).Value = 1080.1546453194
$ws.Cells.Item(7, 18).Value = 9721.3918078746
$ws.Cells.Item(7, 19).Value = 0.1057760128796201
$ws.Cells.Item(7, 20).Value = 0.1057760128796201

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 41.21033366666666
$ws.Cells.Item(8, 8).Value = 123.631001
$ws.Cells.Item(8, 9).Value = 0.1457656157184491
$ws.Cells.Item(8, 10).Value = 0.1457656157184491
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 2.614605666666666
$ws.Cells.Item(8, 14).Value = 7.843817
$ws.Cells.Item(8, 15).Value = 0.07238665237615237
$ws.Cells.Item(8, 16).Value = 0.07238665237615237
$ws.Cells.Item(8, 17).Value = 107.7487719300908
$ws.Cells.Item(8, 18).Value = 969.738947370817
$ws.Cells.Item(8, 19).Value = 0.01055148495340719
$ws.Cells.Item(8, 20).Value = 0.01055148495340719

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 41.21033366666666
$ws.Cells.Item(9, 8).Value = 123.631001
$ws.Cells.Item(9, 9).Value = 0.1457656157184491
$ws.Cells.Item(9, 10).Value = 0.1457656157184491
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 14.70158366666666
$ws.Cells.Item(9, 14).Value = 44.10475099999999
$ws.Cells.Item(9, 15).Value = 0.4070206225838464
$ws.Cells.Item(9, 16).Value = 0.4070206225838464
$ws.Cells.Item(9, 17).Value = 605.85716833175
$ws.Cells.Item(9, 18).Value = 5452.71451498575
$ws.Cells.Item(9, 19).Value = 0.05932961166104087
$ws.Cells.Item(9, 20).Value = 0.05932961166104087

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 41.21033366666666
$ws.Cells.Item(10, 8).Value = 123.631001
$ws.Cells.Item(10, 9).Value = 0.1457656157184491
$ws.Cells.Item(10, 10).Value = 0.1457656157184491
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 18.80380766666667
$ws.Cells.Item(10, 14).Value = 56.411423
$ws.Cells.Item(10, 15).Value = 0.5205927250400012
$ws.Cells.Item(10, 16).Value = 0.5205927250400012
$ws.Cells.Item(10, 17).Value = 774.9111881471581
$ws.Cells.Item(10, 18).Value = 6974.200693324423
$ws.Cells.Item(10, 19).Value = 0.07588451910400106
$ws.Cells.Item(10, 20).Value = 0.07588451910400106
